$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the raw metric values (B2:B13) ---
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsMetrics.Range("B2").Value  = 195204.78
$wsMetrics.Range("B3").Value  = 160540.29
$wsMetrics.Range("B4").Value  = 62035.630000000005
$wsMetrics.Range("B5").Value  = 7713
$wsMetrics.Range("B6").Value  = 4562336.25
$wsMetrics.Range("B7").Value  = 3850358.96
$wsMetrics.Range("B8").Value  = 1332637.7700000003
$wsMetrics.Range("B9").Value  = 176714
$wsMetrics.Range("B10").Value = 33027660.050999828
$wsMetrics.Range("B11").Value = 19880229.030000001
$wsMetrics.Range("B12").Value = 11614346.660000002
$wsMetrics.Range("B13").Value = 1274341

# --- today sheet: drop the stray B2 label cell ---
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("B2").ClearContents()

# --- restore selections / view state to match the captured session ---
$wsToday.Activate()
[void]$wsToday.Range("G13").Select()

$wsMetrics.Activate()
[void]$wsMetrics.Range("H22").Select()
